$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "UI work is done" - mark the remaining finished pages' Status column
# ("Welcome", "Create Community" follow-up row and "User Main Page") as
# Done, copying the look of the existing Done cells (F6 / F21) so the
# new cells get the same number format + borders as the rest of that
# column instead of the blank placeholder style.
function Set-Done([string]$targetAddr, [string]$sourceAddr) {
    $src = $ws.Range($sourceAddr)
    $dst = $ws.Range($targetAddr)
    $dst.Value = "Done"
    $dst.NumberFormat = $src.NumberFormat
    foreach ($edge in 7, 8, 9, 10) {
        $dst.Borders.Item($edge).LineStyle = $src.Borders.Item($edge).LineStyle
        $dst.Borders.Item($edge).Weight = $src.Borders.Item($edge).Weight
        $dst.Borders.Item($edge).Color = $src.Borders.Item($edge).Color
    }
}

Set-Done "F3" "F6"
Set-Done "F9" "F21"
Set-Done "F19" "F21"

# Leave the selection on the last cell that was touched.
$ws.Range("F9").Select()
